$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "16.03.2023"
$ws.Range("C3").Value = "https://gitlab.intra.infineon.com/digital-reference/process_model_version/-/commit/ed806f827269d7f72e29e68f9f504f8efbab3dc6"
$ws.Range("D3").Value = "0e08fd8eee36e999283f9dd25c3209735b91da1642e3b67aadf38a56da7ed5d6"
